{"js": "// Replace the two-digit multiplication equations in the document's table\n// with the newly generated set of equations (same count/order as before).\nconst replacements = [\n  [\"58\u00d754=3132\", \"38\u00d792=3496\"],\n  [\"60\u00d737=2220\", \"67\u00d727=1809\"],\n  [\"99\u00d783=8217\", \"68\u00d799=6732\"],\n  [\"22\u00d734=748\", \"29\u00d737=1073\"],\n  [\"59\u00d719=1121\", \"60\u00d782=4920\"],\n  [\"59\u00d782=4838\", \"16\u00d756=896\"],\n  [\"76\u00d751=3876\", \"82\u00d722=1804\"],\n  [\"52\u00d717=884\", \"65\u00d798=6370\"],\n  [\"90\u00d716=1440\", \"44\u00d783=3652\"],\n  [\"51\u00d717=867\", \"75\u00d725=1875\"],\n  [\"22\u00d798=2156\", \"58\u00d784=4872\"],\n  [\"89\u00d735=3115\", \"78\u00d742=3276\"],\n  [\"22\u00d738=836\", \"55\u00d723=1265\"],\n  [\"99\u00d733=3267\", \"76\u00d718=1368\"],\n  [\"65\u00d774=4810\", \"87\u00d787=7569\"],\n  [\"23\u00d748=1104\", \"23\u00d722=506\"],\n  [\"98\u00d734=3332\", \"36\u00d757=2052\"],\n  [\"57\u00d730=1710\", \"47\u00d788=4136\"],\n  [\"77\u00d759=4543\", \"55\u00d772=3960\"],\n  [\"61\u00d738=2318\", \"56\u00d789=4984\"],\n  [\"89\u00d792=8188\", \"21\u00d747=987\"],\n  [\"11\u00d779=869\", \"53\u00d747=2491\"],\n  [\"50\u00d730=1500\", \"45\u00d718=810\"],\n  [\"72\u00d779=5688\", \"16\u00d733=528\"],\n  [\"69\u00d762=4278\", \"72\u00d729=2088\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication equations in the document's table\n# with the newly generated set of equations (same count/order as before).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"58\u00d754=3132\", \"38\u00d792=3496\"),\n  @(\"60\u00d737=2220\", \"67\u00d727=1809\"),\n  @(\"99\u00d783=8217\", \"68\u00d799=6732\"),\n  @(\"22\u00d734=748\", \"29\u00d737=1073\"),\n  @(\"59\u00d719=1121\", \"60\u00d782=4920\"),\n  @(\"59\u00d782=4838\", \"16\u00d756=896\"),\n  @(\"76\u00d751=3876\", \"82\u00d722=1804\"),\n  @(\"52\u00d717=884\", \"65\u00d798=6370\"),\n  @(\"90\u00d716=1440\", \"44\u00d783=3652\"),\n  @(\"51\u00d717=867\", \"75\u00d725=1875\"),\n  @(\"22\u00d798=2156\", \"58\u00d784=4872\"),\n  @(\"89\u00d735=3115\", \"78\u00d742=3276\"),\n  @(\"22\u00d738=836\", \"55\u00d723=1265\"),\n  @(\"99\u00d733=3267\", \"76\u00d718=1368\"),\n  @(\"65\u00d774=4810\", \"87\u00d787=7569\"),\n  @(\"23\u00d748=1104\", \"23\u00d722=506\"),\n  @(\"98\u00d734=3332\", \"36\u00d757=2052\"),\n  @(\"57\u00d730=1710\", \"47\u00d788=4136\"),\n  @(\"77\u00d759=4543\", \"55\u00d772=3960\"),\n  @(\"61\u00d738=2318\", \"56\u00d789=4984\"),\n  @(\"89\u00d792=8188\", \"21\u00d747=987\"),\n  @(\"11\u00d779=869\", \"53\u00d747=2491\"),\n  @(\"50\u00d730=1500\", \"45\u00d718=810\"),\n  @(\"72\u00d779=5688\", \"16\u00d733=528\"),\n  @(\"69\u00d762=4278\", \"72\u00d729=2088\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $null = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
